{"js": "// Replace the 15 lattice-multiplication problems (5 rows x 3 cols) in the\n// document's single table with a new set of problems, keeping the same\n// \"A x B / <digits of B> / ---- / <digits of A>|    |\" layout/formatting.\n//\n// New problems, addressed by [row, col] (0-based) -> [A, B]:\nconst newProblems = [\n  [0, 0, 67, 36],\n  [0, 1, 15, 28],\n  [0, 2, 24, 19],\n  [1, 0, 72, 75],\n  [1, 1, 38, 55],\n  [1, 2, 94, 81],\n  [2, 0, 13, 88],\n  [2, 1, 23, 92],\n  [2, 2, 23, 16],\n  [3, 0, 37, 23],\n  [3, 1, 90, 38],\n  [3, 2, 57, 30],\n  [4, 0, 69, 83],\n  [4, 1, 98, 33],\n  [4, 2, 43, 84],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, a, b] of newProblems) {\n  const aStr = String(a);\n  const bStr = String(b);\n  const lines = [\n    `${a} x ${b}`,\n    `  ${bStr[0]}    ${bStr[1]}`,\n    `  ----`,\n    `${aStr[0]}|    |`,\n    `${aStr[1]}|    |`,\n  ];\n  // Join lines with a manual line break (\\v == Word's soft line break,\n  // serialized as <w:br/> between <w:t> runs).\n  const newText = lines.join(\"\\v\");\n\n  const cell = table.getCell(row, col);\n  const range = cell.body.getRange();\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems (5 rows x 3 cols) in the\n# document's single table with a new set of problems, keeping the same\n# \"A x B / <digits of B> / ---- / <digits of A>|    |\" layout/formatting.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# New problems, addressed by 1-based (row, col) -> (A, B)\n$newProblems = @(\n  @(1, 1, 67, 36),\n  @(1, 2, 15, 28),\n  @(1, 3, 24, 19),\n  @(2, 1, 72, 75),\n  @(2, 2, 38, 55),\n  @(2, 3, 94, 81),\n  @(3, 1, 13, 88),\n  @(3, 2, 23, 92),\n  @(3, 3, 23, 16),\n  @(4, 1, 37, 23),\n  @(4, 2, 90, 38),\n  @(4, 3, 57, 30),\n  @(5, 1, 69, 83),\n  @(5, 2, 98, 33),\n  @(5, 3, 43, 84)\n)\n\n$vbreak = [char]11\n\nforeach ($p in $newProblems) {\n  $row = $p[0]\n  $col = $p[1]\n  $a = [string]$p[2]\n  $b = [string]$p[3]\n\n  $bFirst = $b.Substring(0,1)\n  $bSecond = $b.Substring(1,1)\n  $aFirst = $a.Substring(0,1)\n  $aSecond = $a.Substring(1,1)\n\n  $line1 = \"$a x $b\"\n  $line2 = \"  $bFirst    $bSecond\"\n  $line3 = \"  ----\"\n  $line4 = \"$aFirst|    |\"\n  $line5 = \"$aSecond|    |\"\n\n  $newText = \"$line1$vbreak$line2$vbreak$line3$vbreak$line4$vbreak$line5\"\n\n  $cell = $tbl.Cell($row, $col)\n  $cell.Range.Text = $newText\n}\n"}
